# edit.ps1 - Applies the "overview -> detail" class-page wording update
# to the UC description document, via Word COM-interop (Find/Replace).
#
# The underlying edit replaces "tổng quan" (overview) wording around the
# class page with "chi tiết" (detail) wording in several use-case tables,
# renames a use case title, reworks a trigger/flow around a "Class Info"
# button (replacing separate "Teacher Info"/"Student Info" buttons), and
# merges two runs that already had identical combined text.

$d = $word.ActiveDocument

# Replace the Nth (1-based) occurrence of $SearchText in the document
# with $ReplaceText. Returns $true if a match was replaced.
function Replace-NthOccurrence {
    param(
        [string]$SearchText,
        [int]$Occurrence,
        [string]$ReplaceText
    )
    $rng = $d.Content
    $rng.Start = 0
    $rng.End = $d.Content.End
    $idx = 0
    $target = $null
    while ($rng.Find.Execute($SearchText)) {
        $idx = $idx + 1
        if ($idx -eq $Occurrence) {
            $target = $d.Range($rng.Start, $rng.End)
            break
        }
        $rng.Start = $rng.End
        $rng.End = $d.Content.End
    }
    if ($target -ne $null) {
        $target.Text = $ReplaceText
        return $true
    }
    return $false
}

# Replace every occurrence of $SearchText in the document with $ReplaceText.
function Replace-AllOccurrences {
    param(
        [string]$SearchText,
        [string]$ReplaceText
    )
    $d.Content.Find.Execute($SearchText, $true, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2) | Out-Null
}

# --- 1. "3. ... xem thông tin tổng quan, thêm học viên ..." -------------
Replace-AllOccurrences `
    "3. Người dùng có thể chọn 1 lớp học bất kỳ để xem thông tin tổng quan, thêm học viên hoặc điểm danh." `
    "3. Người dùng có thể chọn 1 lớp học bất kỳ để xem thông tin chi tiết, thêm học viên hoặc điểm danh."

# --- 2. Use case title: "Xem thông tin giáo viên/học viên" -------------
Replace-AllOccurrences `
    "Xem thông tin giáo viên/học viên" `
    "Xem thông tin chi tiết lớp học"

# --- 3. Merge two runs with already-matching combined text -------------
Replace-AllOccurrences `
    "Người dùng muốn xem thông tin chi tiết của giáo viên/học viên." `
    "Người dùng muốn xem thông tin chi tiết của giáo viên/học viên."

# --- 4. Trigger: Teacher Info/Student Info -> Class Info ----------------
Replace-AllOccurrences `
    "Người dùng bấm vào nút “Teacher Info”/“Student Info” tại trang thông tin tổng quan của lớp học." `
    "Người dùng bấm vào nút “Class Info” để xem thông tin chi tiết lớp học."

# --- 5/8/10/12. Four identical "đang ở tại trang thông tin tổng quan lớp
#     học." sentences that diverge into different replacements -----------
$dangODuocText = "Người dùng đang ở tại trang thông tin tổng quan lớp học."
Replace-NthOccurrence $dangODuocText 1 "Người dùng đang ở tại trang danh sách các lớp mình phụ trách."
Replace-NthOccurrence $dangODuocText 1 "Người dùng đang ở tại trang thông tin chi tiết lớp học."
Replace-NthOccurrence $dangODuocText 1 "Người dùng đang ở tại trang thông tin chi tiết lớp học."
Replace-NthOccurrence $dangODuocText 1 "Người dùng đang ở tại trang thông tin chi tiết lớp học."

# --- 6b/6c. Normal flow steps 2 and 3 of "Xem thông tin chi tiết lớp học"
Replace-AllOccurrences `
    "2. Người dùng chọn nút “Teacher Info” hoặc“Student Info”" `
    "2. Người dùng chọn nút “Class Info”."

Replace-AllOccurrences `
    "3. Hệ thống hiện 1 trang mini thông tin của giáo viên/học viên." `
    "3. Hệ thống chuyển sang trang hiển thị thông tin chi tiết của một lớp."

# --- 6a/9. Two identical "đang ở trang hiển thị thông tin tổng quan của
#     lớp học." sentences that diverge into different replacements -------
$hienThiText = "1. Người dùng đang ở trang hiển thị thông tin tổng quan của lớp học."
Replace-NthOccurrence $hienThiText 1 "1. Người dùng đang ở trang hiển thị danh sách lớp học."
Replace-NthOccurrence $hienThiText 1 "1. Người dùng đang ở trang hiển thị thông tin chi tiết của lớp học."

# --- 7. "Add a student" trigger -----------------------------------------
Replace-AllOccurrences `
    "Người dùng bấm vào nút “Add a student” tại trang thông tin tổng quan của lớp học." `
    "Người dùng bấm vào nút “Add a student” tại trang thông tin chi tiết của lớp học."

# --- 11/13. Attendance "Check attendace" normal-flow step 1 (both
#     occurrences use the exact same replacement) ------------------------
Replace-AllOccurrences `
    "1. Tại trang thông tin tổng quan lớp học, người dùng xác định buổi học cụ thể và bấm vào nút “Check attendace” ở cuối mỗi hàng của danh sách buổi học." `
    "1. Tại trang thông tin chi tiết lớp học, người dùng xác định buổi học cụ thể và bấm vào nút “Check attendace” ở cuối mỗi hàng của danh sách buổi học."

Write-Output "Done."
